$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cancellation record shown on row 2 with the latest transaction
# (cuenta, Transaccion and Fecha updated; Estado stays PASSED)
$ws.Range("C2").Value = "'1010826124"
$ws.Range("G2").Value = "AAACT23195WF05QJC"
$ws.Range("H2").Value = "14 jul. 2023, 09:30:21"

# Move/leave the active selection on E9
$ws.Range("E9").Select()
